# Add a second worksheet that holds a condensed (first 6 periods only,
# summary rows only) copy of Sheet1's seasonality table.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Duplicate Sheet1 (this carries over all formatting/styles/data) and place
# the copy immediately after Sheet1.
$sheet1.Copy($null, $sheet1)
$sheet2 = $wb.Worksheets.Item(2)
$sheet2.Name = "Sheet2"

# Trim the copy down to the condensed table: keep only Order 1 - Order 6
# (rows 2:7) plus the summary rows (Demand/Period/Setup_Cost/Holding_Cost),
# and only the first six period columns (A:G).
$sheet2.Range("A8:A13").EntireRow.Delete()
$sheet2.Range("H1:M1").EntireColumn.Delete()

# Restore Sheet1's selection (no longer the active tab / active cell).
$sheet1.Range("A1:G17").Select() | Out-Null

# Sheet2 becomes the active tab with its own fresh selection.
$sheet2.Activate() | Out-Null
$sheet2.Range("B19").Select() | Out-Null
